# Results for NTRANS konferanse - Inkludert VSS
# Adds a new "EV_scenario" worksheet (between "three_scenarios_new" and
# "variability") with a single base/base/base scenario row, and nudges the
# remembered selection on a couple of existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "EV_scenario" sheet right after "three_scenarios_new"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("three_scenarios_new")
$evSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$evSheet.Name = "EV_scenario"

# Header row - same column layout as the other scenario sheets.
$evSheet.Range("A1").Value = "Scenario"
$evSheet.Range("B1").Value = "Name"
$evSheet.Range("C1").Value = "Probability"
$evSheet.Range("D1").Value = "Cost_Established"
$evSheet.Range("E1").Value = "Cost_Battery"
$evSheet.Range("F1").Value = "Cost_Hydrogen"
$evSheet.Range("G1").Value = "Cost_Biofuel"
$evSheet.Range("H1").Value = "Maturity_Established"
$evSheet.Range("I1").Value = "Maturity_Battery"
$evSheet.Range("J1").Value = "Maturity_Hydrogen"
$evSheet.Range("K1").Value = "Maturity_Biofuel"
$evSheet.Range("A1:K1").Font.Bold = $true

# Single data row: scenario 1, "MMM", probability 1, all costs 1, all
# maturities "base".
$evSheet.Range("A2").Value = 1
$evSheet.Range("B2").Value = "MMM"
$evSheet.Range("C2").Value = 1
$evSheet.Range("D2").Value = 1
$evSheet.Range("E2").Value = 1
$evSheet.Range("F2").Value = 1
$evSheet.Range("G2").Value = 1
$evSheet.Range("H2").Value = "base"
$evSheet.Range("I2").Value = "base"
$evSheet.Range("J2").Value = "base"
$evSheet.Range("K2").Value = "base"

# ---------------------------------------------------------------------
# 2. Update remembered selections on a couple of sheets
# ---------------------------------------------------------------------
$baseSheet = $wb.Worksheets.Item("scenarios_base")
$baseSheet.Activate() | Out-Null
$baseSheet.Range("A15:K15").Select() | Out-Null

$evSheet.Activate() | Out-Null
$evSheet.Range("D11").Select() | Out-Null

$newScenarioSheet = $wb.Worksheets.Item("three_scenarios_new")
$newScenarioSheet.Activate() | Out-Null
$newScenarioSheet.Range("G9").Select() | Out-Null

Write-Output "EV_scenario sheet added"
